# Add data for 2025-11-23
# Updates column L (year 2025) violent-crime running totals across the
# citywide summary, by-neighborhood summary, and each neighborhood detail sheet.
# Two cells (B4/B7 on "Citywide Totals" and B63/B101 on "By Neighborhood") are
# historical corrections to the 2015 column bundled with this data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 5980
$ws.Range("L3").Value = 6521
$ws.Range("B4").Value = 1720
$ws.Range("L4").Value = 1599
$ws.Range("L5").Value = 388
$ws.Range("L6").Value = 5356
$ws.Range("B7").Value = 23352
$ws.Range("L7").Value = 19844

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L6").Value = 154
$ws.Range("L7").Value = 642
$ws.Range("L8").Value = 1314
$ws.Range("L11").Value = 333
$ws.Range("L19").Value = 537
$ws.Range("L20").Value = 504
$ws.Range("L29").Value = 1117
$ws.Range("L31").Value = 196
$ws.Range("L33").Value = 894
$ws.Range("L34").Value = 112
$ws.Range("L36").Value = 251
$ws.Range("L37").Value = 751
$ws.Range("L42").Value = 635
$ws.Range("L43").Value = 150
$ws.Range("L47").Value = 138
$ws.Range("L48").Value = 260
$ws.Range("L49").Value = 106
$ws.Range("L52").Value = 421
$ws.Range("L53").Value = 219
$ws.Range("L54").Value = 433
$ws.Range("L55").Value = 206
$ws.Range("L56").Value = 19
$ws.Range("B63").Value = 424
$ws.Range("L65").Value = 389
$ws.Range("L67").Value = 683
$ws.Range("L76").Value = 302
$ws.Range("L78").Value = 257
$ws.Range("L79").Value = 551
$ws.Range("L82").Value = 28
$ws.Range("L85").Value = 982
$ws.Range("L88").Value = 210
$ws.Range("L89").Value = 276
$ws.Range("L90").Value = 206
$ws.Range("L95").Value = 283
$ws.Range("L98").Value = 105
$ws.Range("B101").Value = 23352
$ws.Range("L101").Value = 19844

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L3").Value = 206
$ws.Range("L7").Value = 642

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L2").Value = 125
$ws.Range("L3").Value = 97
$ws.Range("L6").Value = 83
$ws.Range("L7").Value = 333

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L6").Value = 78
$ws.Range("L7").Value = 276

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 294
$ws.Range("L3").Value = 407
$ws.Range("L7").Value = 982

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L2").Value = 130
$ws.Range("L3").Value = 132
$ws.Range("L6").Value = 120
$ws.Range("L7").Value = 421

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L4").Value = 23
$ws.Range("L7").Value = 219

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 394
$ws.Range("L3").Value = 464
$ws.Range("L5").Value = 44
$ws.Range("L6").Value = 323
$ws.Range("L7").Value = 1314

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 240
$ws.Range("L3").Value = 313
$ws.Range("L6").Value = 260
$ws.Range("L7").Value = 894

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("L3").Value = 90
$ws.Range("L6").Value = 67
$ws.Range("L7").Value = 283

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L3").Value = 265
$ws.Range("L6").Value = 200
$ws.Range("L7").Value = 751

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L2").Value = 142
$ws.Range("L7").Value = 389

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L2").Value = 80
$ws.Range("L7").Value = 196

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L3").Value = 267
$ws.Range("L7").Value = 683

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("L2").Value = 35
$ws.Range("L7").Value = 106

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L3").Value = 110
$ws.Range("L7").Value = 433

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 330
$ws.Range("L3").Value = 434
$ws.Range("L7").Value = 1117

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L3").Value = 68
$ws.Range("L7").Value = 260

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L4").Value = 27
$ws.Range("L7").Value = 537

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L3").Value = 60
$ws.Range("L7").Value = 302

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("L2").Value = 60
$ws.Range("L7").Value = 154

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L3").Value = 219
$ws.Range("L7").Value = 635

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L2").Value = 68
$ws.Range("L4").Value = 29
$ws.Range("L7").Value = 257

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L4").Value = 19
$ws.Range("L7").Value = 206

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L5").Value = 17
$ws.Range("L7").Value = 551

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L3").Value = 176
$ws.Range("L7").Value = 504

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L3").Value = 81
$ws.Range("L6").Value = 60
$ws.Range("L7").Value = 251

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("L2").Value = 39
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 112

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("L2").Value = 52
$ws.Range("L7").Value = 138

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 105

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L2").Value = 65
$ws.Range("L7").Value = 210

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("L3").Value = 59
$ws.Range("L7").Value = 206

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("L3").Value = 51
$ws.Range("L4").Value = 23
$ws.Range("L7").Value = 150

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("L2").Value = 4
$ws.Range("L7").Value = 28

$ws = $wb.Worksheets.Item('Magnificent Mile')
$ws.Range("L4").Value = 3
$ws.Range("L7").Value = 19
